# Weekly update: prepend 2 new price rows (most recent week) above the
# existing data in the "Fruta, Agrícola del Norte S.A. de Arica - Limón"
# sheet, pushing all prior rows down by two. Dimension grows from
# A1:T123 to A1:T125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the first data row (row 29), shifting the
# existing rows 29-123 down to 31-125.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# --- New row 29 ---------------------------------------------------------
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C29").Value = 'Arica y Parinacota'
$ws.Range("D29").Value = 44414
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 'Fruta'
$ws.Range("G29").Value = 100102
$ws.Range("H29").Value = 'Cítricos'
$ws.Range("I29").Value = 100102003
$ws.Range("J29").Value = 'Limón'
$ws.Range("K29").Value = 'Sutil De Gase'
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 250
$ws.Range("N29").Value = 31000
$ws.Range("O29").Value = 32000
$ws.Range("P29").Value = 31500
$ws.Range("Q29").Value = '$/caja 24 kilos'
$ws.Range("R29").Value = 'Perú'
$ws.Range("S29").Value = 1312
$ws.Range("T29").Value = 24

# --- New row 30 ---------------------------------------------------------
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C30").Value = 'Arica y Parinacota'
$ws.Range("D30").Value = 44414
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = 'Fruta'
$ws.Range("G30").Value = 100102
$ws.Range("H30").Value = 'Cítricos'
$ws.Range("I30").Value = 100102003
$ws.Range("J30").Value = 'Limón'
$ws.Range("K30").Value = 'Tahití'
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 32000
$ws.Range("O30").Value = 33000
$ws.Range("P30").Value = 32500
$ws.Range("Q30").Value = '$/caja 24 kilos'
$ws.Range("R30").Value = 'Perú'
$ws.Range("S30").Value = 1354
$ws.Range("T30").Value = 24
